$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 922
$ws.Range("E2").Value = -94
$ws.Range("F2").Value = -94
$ws.Range("G2").Value = -75
$ws.Range("H2").Value = -65
$ws.Range("I2").Value = -62
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 8745
$ws.Range("L2").Value = 7316
$ws.Range("M2").Value = 1429
$ws.Range("N2").Value = 1410
$ws.Range("O2").Value = 19
$ws.Range("P2").Value = 650
$ws.Range("Q2").Value = -608
$ws.Range("R2").Value = -548
$ws.Range("S2").Value = 902
$ws.Range("T2").Value = 188
$ws.Range("V2").Value = 1261
$ws.Range("W2").Value = -10.18
$ws.Range("X2").Value = -7.01
$ws.Range("Y2").Value = -4.4
$ws.Range("Z2").Value = -0.87
$ws.Range("AA2").Value = 512.09
$ws.Range("AB2").Value = 124.95
$ws.Range("AC2").Value = -97
$ws.Range("AD2").Value = -7.61
$ws.Range("AE2").Value = 2286
$ws.Range("AF2").Value = 0.32
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 63737067

# Row 3
$ws.Range("D3").Value = 898
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 45
$ws.Range("H3").Value = 31
$ws.Range("I3").Value = 32
$ws.Range("J3").Value = -1
$ws.Range("K3").Value = 7029
$ws.Range("L3").Value = 5584
$ws.Range("M3").Value = 1445
$ws.Range("N3").Value = 1428
$ws.Range("O3").Value = 17
$ws.Range("P3").Value = 650
$ws.Range("Q3").Value = -442
$ws.Range("R3").Value = -5
$ws.Range("S3").Value = 500
$ws.Range("T3").Value = 2
$ws.Range("V3").Value = 692
$ws.Range("W3").Value = 0.75
$ws.Range("X3").Value = 3.47
$ws.Range("Y3").Value = 2.27
$ws.Range("Z3").Value = 0.41
$ws.Range("AA3").Value = 386.43
$ws.Range("AB3").Value = 127.45
$ws.Range("AC3").Value = 51
$ws.Range("AD3").Value = 17.22
$ws.Range("AE3").Value = 2316
$ws.Range("AF3").Value = 0.38
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 63737067

# Row 4
$ws.Range("D4").Value = 636
$ws.Range("E4").Value = -41
$ws.Range("F4").Value = -41
$ws.Range("G4").Value = -46
$ws.Range("H4").Value = -43
$ws.Range("I4").Value = -42
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 5598
$ws.Range("L4").Value = 4097
$ws.Range("M4").Value = 1501
$ws.Range("N4").Value = 1415
$ws.Range("O4").Value = 86
$ws.Range("P4").Value = 650
$ws.Range("Q4").Value = 85
$ws.Range("R4").Value = 48
$ws.Range("S4").Value = -109
$ws.Range("T4").Value = 2
$ws.Range("V4").Value = 1274
$ws.Range("W4").Value = -6.38
$ws.Range("X4").Value = -6.78
$ws.Range("Y4").Value = -2.97
$ws.Range("Z4").Value = -0.67
$ws.Range("AA4").Value = 273
$ws.Range("AB4").Value = 136.04
$ws.Range("AC4").Value = -66
$ws.Range("AD4").Value = -11.06
$ws.Range("AE4").Value = 2295
$ws.Range("AF4").Value = 0.32
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 63737067

# Row 5
$ws.Range("D5").Value = 454
$ws.Range("E5").Value = -92
$ws.Range("F5").Value = -92
$ws.Range("G5").Value = -48
$ws.Range("H5").Value = -42
$ws.Range("I5").Value = -46
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 8565
$ws.Range("L5").Value = 7042
$ws.Range("M5").Value = 1522
$ws.Range("N5").Value = 1387
$ws.Range("O5").Value = 135
$ws.Range("P5").Value = 650
$ws.Range("Q5").Value = -1131
$ws.Range("R5").Value = 92
$ws.Range("S5").Value = 1287
$ws.Range("T5").Value = 3
$ws.Range("V5").Value = 1127
$ws.Range("W5").Value = -20.32
$ws.Range("X5").Value = -9.34
$ws.Range("Y5").Value = -3.28
$ws.Range("Z5").Value = -0.65
$ws.Range("AA5").Value = 462.6
$ws.Range("AB5").Value = 139.34
$ws.Range("AC5").Value = -78
$ws.Range("AD5").Value = -14.74
$ws.Range("AE5").Value = 2853
$ws.Range("AF5").Value = 0.4
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 50693589

# Row 6
$ws.Range("D6").Value = 343
$ws.Range("E6").Value = -114
$ws.Range("F6").Value = -114
$ws.Range("G6").Value = -108
$ws.Range("H6").Value = -94
$ws.Range("I6").Value = -99
$ws.Range("K6").Value = 4998
$ws.Range("L6").Value = 3923
$ws.Range("M6").Value = 1075
$ws.Range("N6").Value = 986
$ws.Range("P6").Value = 520
$ws.Range("Q6").Value = 1491
$ws.Range("R6").Value = 7
$ws.Range("S6").Value = -1705
$ws.Range("V6").Value = 1092
$ws.Range("W6").Value = -33.36
$ws.Range("X6").Value = -27.4
$ws.Range("Y6").Value = -8.32
$ws.Range("Z6").Value = -1.46
$ws.Range("AA6").Value = 364.78
$ws.Range("AB6").Value = 113.46
$ws.Range("AC6").Value = -195
$ws.Range("AD6").Value = -9.35
$ws.Range("AE6").Value = 2027
$ws.Range("AF6").Value = 0.9
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 50693589

# Clear cells removed from rows 2-6
$ws.Range("U2").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("T6:U6").ClearContents()

# Rows 7,8,9: clear all data cells (D:AJ), keep A,B,C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()